$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at position 15-16, shifting existing rows 15-78 down to 17-80
$ws.Rows("15:16").Insert()

# Fill the two newly inserted rows with the new weekly price entries
$row15 = New-Object 'object[,]' 1,18
$row15[0,0] = 11
$row15[0,1] = "Vega Monumental Concepción"
$row15[0,2] = "Bíobío"
$row15[0,3] = 44881
$row15[0,4] = 8
$row15[0,5] = 100112037
$row15[0,6] = "Cebollín"
$row15[0,7] = "Sin especificar"
$row15[0,8] = "Primera"
$row15[0,9] = 200
$row15[0,10] = 700
$row15[0,11] = 800
$row15[0,12] = 750
$row15[0,13] = "`$/paquete 6 unidades"
$row15[0,14] = "Región de Ñuble"
$row15[0,15] = 125
$row15[0,16] = 6
$row15[0,17] = "Hortaliza"
$ws.Range("A15:R15").Value = $row15

$row16 = New-Object 'object[,]' 1,18
$row16[0,0] = 11
$row16[0,1] = "Vega Monumental Concepción"
$row16[0,2] = "Bíobío"
$row16[0,3] = 44881
$row16[0,4] = 8
$row16[0,5] = 100112037
$row16[0,6] = "Cebollín"
$row16[0,7] = "Sin especificar"
$row16[0,8] = "Segunda"
$row16[0,9] = 100
$row16[0,10] = 600
$row16[0,11] = 600
$row16[0,12] = 600
$row16[0,13] = "`$/paquete 6 unidades"
$row16[0,14] = "Región de Ñuble"
$row16[0,15] = 100
$row16[0,16] = 6
$row16[0,17] = "Hortaliza"
$ws.Range("A16:R16").Value = $row16

